$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("plotConfiguration")

# Replace the four separate limit columns (xLimLower, xLimUpper, yLimLower,
# yLimUpper) with two combined columns (xAxisLimits, yAxisLimits) that hold
# comma-separated "lower, upper" values.
$ws.Range("K1").Value = ""
$ws.Range("L1").Value = ""
$ws.Range("I1").Value = "xAxisLimits"
$ws.Range("J1").Value = "yAxisLimits"
$ws.Range("I2").Value = "0, 24"

# Make plotConfiguration the active sheet/tab instead of DataCombined.
$ws.Activate()
$ws.Range("I3").Select()
